# Weekly update: add a new date column (AH) with the 22_06_2021 admissions data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new week.
$ws.Range("AH1").Value = "22_06_2021"

# Per age-group counts for the new week.
$ws.Range("AH2").Value = 13
$ws.Range("AH3").Value = 19
$ws.Range("AH4").Value = 37
$ws.Range("AH5").Value = 56
$ws.Range("AH6").Value = 137
$ws.Range("AH7").Value = 300
$ws.Range("AH8").Value = 426
$ws.Range("AH9").Value = 578
$ws.Range("AH10").Value = 180
$ws.Range("AH11").Value = 13

# Total row: sum of the new column.
$ws.Range("AH12").Formula = "=SUM(AH2:AH11)"

$excel.Calculate()

# Move the current selection to match where the editor ended up after
# entering the new week's data.
$ws.Range("AG14").Select() | Out-Null
